# Reorganización completa: limpieza de módulos antiguos, nuevas entregas y optimización
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Datos" to "potreros"
$ws.Name = "potreros"

$headerRange = $ws.Range("A1:J1")

# Drop the old header styling (bold white font on blue fill, centered) --
# back to the sheet's plain default formatting.
$headerRange.ClearFormats()

# Remove the explicit 20-wide column sizing entirely (delete + reinsert
# brings the columns back to the workbook's default width with no leftover
# per-column override).
$headerRange.EntireColumn.Delete()
$ws.Range("A1:J1").Insert()

# New plain (snake_case, lowercase, unaccented) header labels replacing the
# old capitalised / accented ones
$headers = @("codigo", "finca", "nombre", "sector", "area_hectareas", "capacidad_maxima", "tipo_pasto", "descripcion", "estado", "comentario")

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}
